$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the quantity/amount figures for the exam-work bill that were left
# blank in the template. Entering these values also feeds the existing
# formulas in column I (and the grand-total in I32), which recalculate
# automatically.
$ws.Range("G9").Value = 117
$ws.Range("G12").Value = 117
$ws.Range("G16").Value = 27
$ws.Range("G17").Value = 58.5
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
$ws.Range("G29").Value = 10
